# Update the worksheet date and the 25 division problems.
$d = $word.ActiveDocument

# 1. Update the date heading.
$d.Content.Find.Execute("2024-02-20 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-02-21 Wednesday", 2)

# 2. Update the division problems in the table, cell by cell (position-based,
#    so the unique-text-per-cell assumption doesn't matter and there is no
#    risk of one replacement's output colliding with another's input text).
$t = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)

$values = @(
    @("93÷6=", "45÷4=", "11÷2=", "45÷8=", "36÷8="),
    @("96÷7=", "43÷6=", "83÷9=", "93÷2=", "26÷4="),
    @("63÷4=", "36÷2=", "57÷6=", "50÷7=", "70÷8="),
    @("26÷7=", "92÷5=", "15÷3=", "39÷4=", "62÷9="),
    @("46÷8=", "91÷5=", "33÷4=", "11÷4=", "27÷2=")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $tableRow = $rows[$i]
    $rowValues = $values[$i]
    for ($col = 1; $col -le 5; $col++) {
        $t.Cell($tableRow, $col).Range.Text = $rowValues[$col - 1]
    }
}
